$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logindata")

# New row 7 mirrors row 6 (same B value/style, same C "Y"), but with a new A value.
$ws.Range("A7").Value = "test_edit_user"

$ws.Range("B7").Value = "TestUsr1,Mitel@123,Mitel@gmail.com,MitelFirst,MitelSecond,`nBusiness1,Business2,+917975935256,+918105855417"
$ws.Range("B7").WrapText = $true

$ws.Range("C7").Value = "Y"

$ws.Rows.Item(7).RowHeight = 58

$ws.Range("F7").Select()
